$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 125000184
$ws.Range("I6").Value = 125000184
$ws.Range("K6").Value = 375000552
$ws.Range("M6").Value = -375000440
$ws.Range("H39").Value = 100.5
$ws.Range("I39").Value = 100.5
$ws.Range("K39").Value = 301.5
$ws.Range("M39").Value = -5.5
$ws.Range("H40").Value = 2263.5833
$ws.Range("I40").Value = 1280.4286
$ws.Range("K40").Value = 1280.4286
$ws.Range("M40").Value = -1105.4286
$ws.Range("H86").Value = 3877.8
$ws.Range("I86").Value = 3870.1177
$ws.Range("J86").Value = 3921.3333
$ws.Range("K86").Value = 3870.1177
$ws.Range("L86").Value = 3921.3333
$ws.Range("M86").Value = -2747.1177
$ws.Range("N86").Value = -6167.3333
$ws.Range("H89").Value = 3877.8
$ws.Range("I89").Value = 3870.1177
$ws.Range("J89").Value = 3921.3333
$ws.Range("K89").Value = 19350.5885
$ws.Range("L89").Value = 19606.6665
$ws.Range("M89").Value = -13734.5885
$ws.Range("N89").Value = -30838.6665
$ws.Range("H106").Value = 4309.1816
$ws.Range("I106").Value = 4279.6523
$ws.Range("J106").Value = 4377.1
$ws.Range("K106").Value = 4279.6523
$ws.Range("L106").Value = 4377.1
$ws.Range("M106").Value = -3648.6523
$ws.Range("N106").Value = -5639.1
$ws.Range("H127").Value = 1387.9166
$ws.Range("J127").Value = 1497.5
$ws.Range("L127").Value = 4492.5
$ws.Range("N127").Value = -14412.5
$ws.Range("H132").Value = 1756.6875
$ws.Range("I132").Value = 1062.7142
$ws.Range("K132").Value = 3188.1426
$ws.Range("M132").Value = -658.1425999999997
$ws.Range("H137").Value = 1531
$ws.Range("I137").Value = 1376.65
$ws.Range("K137").Value = 4129.950000000001
$ws.Range("M137").Value = -1579.950000000001

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 418.375
$ws.Range("I5").Value = 475
$ws.Range("J5").Value = 399.5
$ws.Range("K5").Value = 475
$ws.Range("L5").Value = 399.5
$ws.Range("M5").Value = -363
$ws.Range("N5").Value = -623.5
$ws.Range("H32").Value = 4698.1035
$ws.Range("I32").Value = 4516.16
$ws.Range("K32").Value = 4516.16
$ws.Range("M32").Value = -4229.16
$ws.Range("H61").Value = 5091.647
$ws.Range("I61").Value = 5091.647
$ws.Range("K61").Value = 5091.647
$ws.Range("M61").Value = -4879.647
$ws.Range("H132").Value = 1483.258
$ws.Range("I132").Value = 1483.258
$ws.Range("K132").Value = 4449.774
$ws.Range("M132").Value = -1919.774
$ws.Range("H136").Value = 5091.647
$ws.Range("I136").Value = 5091.647
$ws.Range("K136").Value = 15274.941
$ws.Range("M136").Value = -12724.941

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 418.375
$ws.Range("I4").Value = 475
$ws.Range("J4").Value = 399.5
$ws.Range("K4").Value = 475
$ws.Range("L4").Value = 399.5
$ws.Range("M4").Value = -360
$ws.Range("N4").Value = -629.5
$ws.Range("H22").Value = 2320
$ws.Range("J22").Value = 2002
$ws.Range("L22").Value = 2002
$ws.Range("N22").Value = -2348
$ws.Range("H35").Value = 57037.2
$ws.Range("J35").Value = 57037.2
$ws.Range("L35").Value = 57037.2
$ws.Range("N35").Value = -57657.2
$ws.Range("H99").Value = 3033.65
$ws.Range("I99").Value = 1309.75
$ws.Range("K99").Value = 1309.75
$ws.Range("M99").Value = 188.25
$ws.Range("H105").Value = 2453.9333
$ws.Range("I105").Value = 2453.9333
$ws.Range("K105").Value = 2453.9333
$ws.Range("M105").Value = -706.9333000000001

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 22113.666
$ws.Range("J38").Value = 20717.572
$ws.Range("L38").Value = 20717.572
$ws.Range("N38").Value = -21471.572
$ws.Range("H46").Value = 22113.666
$ws.Range("J46").Value = 20717.572
$ws.Range("L46").Value = 20717.572
$ws.Range("N46").Value = -21139.572

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7308437.5
$ws.Range("I4").Value = 1622525.5
$ws.Range("J4").Value = 39291692
$ws.Range("K4").Value = 4867576.5
$ws.Range("L4").Value = 117875076
$ws.Range("M4").Value = -4867464.5
$ws.Range("N4").Value = -117875300
$ws.Range("H80").Value = 1639.4
$ws.Range("I80").Value = 1299.25
$ws.Range("K80").Value = 3897.75
$ws.Range("M80").Value = -2961.75
$ws.Range("H83").Value = 1639.4
$ws.Range("I83").Value = 1299.25
$ws.Range("K83").Value = 11693.25
$ws.Range("M83").Value = -7013.25
$ws.Range("H95").Value = 67500
$ws.Range("J95").Value = 67500
$ws.Range("L95").Value = 202500
$ws.Range("N95").Value = -206618
$ws.Range("H113").Value = 1476.8846
$ws.Range("J113").Value = 1740.35
$ws.Range("L113").Value = 5221.049999999999
$ws.Range("N113").Value = -9561.049999999999
$ws.Range("H121").Value = 790.1539
$ws.Range("I121").Value = 475
$ws.Range("J121").Value = 1499.25
$ws.Range("K121").Value = 1425
$ws.Range("L121").Value = 4497.75
$ws.Range("M121").Value = -115
$ws.Range("N121").Value = -7117.75

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 204.11539
$ws.Range("I2").Value = 36.5
$ws.Range("J2").Value = 762.8333
$ws.Range("K2").Value = 36.5
$ws.Range("L2").Value = 762.8333
$ws.Range("M2").Value = 76.5
$ws.Range("N2").Value = -988.8333
$ws.Range("H15").Value = 23640.666
$ws.Range("J15").Value = 23640.666
$ws.Range("L15").Value = 23640.666
$ws.Range("N15").Value = -24216.666
$ws.Range("H20").Value = 1438569.9
$ws.Range("I20").Value = 2506749.8
$ws.Range("J20").Value = 14330
$ws.Range("K20").Value = 2506749.8
$ws.Range("L20").Value = 14330
$ws.Range("M20").Value = -2506504.8
$ws.Range("N20").Value = -14820
$ws.Range("H80").Value = 3836
$ws.Range("I80").Value = 4074.8572
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 4074.8572
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -3076.8572
$ws.Range("N80").Value = -4996
$ws.Range("H81").Value = 23640.666
$ws.Range("J81").Value = 23640.666
$ws.Range("L81").Value = 23640.666
$ws.Range("N81").Value = -25636.666
$ws.Range("H83").Value = 3836
$ws.Range("I83").Value = 4074.8572
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 20374.286
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -15382.286
$ws.Range("N83").Value = -24984
$ws.Range("H84").Value = 23640.666
$ws.Range("J84").Value = 23640.666
$ws.Range("L84").Value = 70921.99800000001
$ws.Range("N84").Value = -80905.99800000001
$ws.Range("H97").Value = 387
$ws.Range("I97").Value = 350
$ws.Range("J97").Value = 461
$ws.Range("K97").Value = 350
$ws.Range("L97").Value = 461
$ws.Range("M97").Value = 146
$ws.Range("N97").Value = -1453
$ws.Range("H132").Value = 2548.1875
$ws.Range("I132").Value = 2548.1875
$ws.Range("K132").Value = 7644.5625
$ws.Range("M132").Value = -5114.5625

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3816.3103
$ws.Range("I40").Value = 2772.3157
$ws.Range("K40").Value = 2772.3157
$ws.Range("M40").Value = -2636.3157
$ws.Range("H46").Value = 1948.8889
$ws.Range("I46").Value = 1359.4
$ws.Range("K46").Value = 1359.4
$ws.Range("M46").Value = -1171.4
$ws.Range("H55").Value = 1739.0588
$ws.Range("I55").Value = 257
$ws.Range("J55").Value = 3856.2856
$ws.Range("K55").Value = 257
$ws.Range("L55").Value = 3856.2856
$ws.Range("M55").Value = -84
$ws.Range("N55").Value = -4202.2856
$ws.Range("H93").Value = 7473.467
$ws.Range("I93").Value = 5000
$ws.Range("J93").Value = 7650.143
$ws.Range("K93").Value = 5000
$ws.Range("L93").Value = 7650.143
$ws.Range("M93").Value = -3752
$ws.Range("N93").Value = -10146.143
$ws.Range("H122").Value = 9612.75
$ws.Range("I122").Value = 9650.333000000001
$ws.Range("K122").Value = 28950.999
$ws.Range("M122").Value = -26500.999
$ws.Range("H136").Value = 17255.242
$ws.Range("I136").Value = 1545.6842
$ws.Range("K136").Value = 4637.0526
$ws.Range("M136").Value = -2087.0526

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4703.7915
$ws.Range("I122").Value = 2725.8
$ws.Range("K122").Value = 8177.400000000001
$ws.Range("M122").Value = -5727.400000000001
